$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("utleielokaler-i-gjesdal-kommune")

# Remove the trailing space from the "Stamphuset " label in column E (row 7)
$ws.Range("E7").Value = "Stamphuset"

# Update the active selection to match the edited cell
$ws.Range("E7").Select()
